# Daily attendance processing - 2025-10-06 10:40:17
# Applies the attendance-session-analysis update:
#  - "Recorded By" string reorders ("system, System" -> "System, system")
#  - Newly-recorded sessions (rows 16, 43, 70) flip from Pending -> Recorded
#    with attendance figures and recorder e-mail filled in (style follows
#    the same visual treatment as the other "Recorded" rows)
#  - A few session rows get a final attendance count top-up
#    (partial counts -> full counts)
#  - Downstream class-statistics table (K:S, rows 6-20) recomputed to match

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) "Recorded By" text - swap order of "system, System" tokens
# ---------------------------------------------------------------------
$ws.Cells.Item(2, 7).Value = "backup@backdoor.com, System, system"
$ws.Cells.Item(29, 7).Value = "backup@backdoor.com, System, system"
$ws.Cells.Item(56, 7).Value = "backup@backdoor.com, System, system"

# "backup@backdoor.com" -> "backup@backdoor.com, System"
$ws.Cells.Item(84, 7).Value = "backup@backdoor.com, System"
$ws.Cells.Item(85, 7).Value = "backup@backdoor.com, System"
$ws.Cells.Item(110, 7).Value = "backup@backdoor.com, System"
$ws.Cells.Item(111, 7).Value = "backup@backdoor.com, System"
$ws.Cells.Item(136, 7).Value = "backup@backdoor.com, System"
$ws.Cells.Item(137, 7).Value = "backup@backdoor.com, System"

# ---------------------------------------------------------------------
# 2) Attendance-count top-ups on already-recorded sessions
# ---------------------------------------------------------------------
$ws.Cells.Item(3, 8).Value = "53/53"
$ws.Cells.Item(56, 8).Value = "55/55"
$ws.Cells.Item(84, 8).Value = "56/56"
$ws.Cells.Item(111, 8).Value = "55/55"

# ---------------------------------------------------------------------
# 3) Sessions that flip from "Pending" to "Recorded" (rows 16, 43, 70)
#    Copy the visual format (fill/font) of an already-"Recorded" row
#    onto each of these rows, then fill in the recorder / count / status.
# ---------------------------------------------------------------------
$ws.Range("A2:I2").Copy()
$ws.Range("A16:I16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(16, 7).Value = "dnasr281@gmail.com"
$ws.Cells.Item(16, 8).Value = "42/53"
$ws.Cells.Item(16, 9).Value = "Recorded"

$ws.Range("A2:I2").Copy()
$ws.Range("A43:I43").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(43, 7).Value = "dnasr281@gmail.com"
$ws.Cells.Item(43, 8).Value = "43/56"
$ws.Cells.Item(43, 9).Value = "Recorded"

$ws.Range("A2:I2").Copy()
$ws.Range("A70:I70").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(70, 7).Value = "dnasr281@gmail.com"
$ws.Cells.Item(70, 8).Value = "36/55"
$ws.Cells.Item(70, 9).Value = "Recorded"

# ---------------------------------------------------------------------
# 4) Class-statistics summary block (columns K:S) recomputed
# ---------------------------------------------------------------------

# Headline counters (K4:L10)
$ws.Cells.Item(6, 12).Value = 80
$ws.Cells.Item(8, 12).Value = 78

# Percentage cells are stored as literal text (e.g. "50.3%"), not numeric
# percentages - use a leading apostrophe to force text entry (Excel's
# standard "treat as text" convention); the quote-prefix style this leaves
# behind on the cell gets cleaned up below.
$ws.Cells.Item(9, 12).Value = "'50.3%"
$ws.Cells.Item(10, 12).Value = "'69.5%"

$ws.Cells.Item(15, 18).Value = "'55.6%"
$ws.Cells.Item(15, 19).Value = "'69.7%"

$ws.Cells.Item(16, 18).Value = "'51.9%"
$ws.Cells.Item(16, 19).Value = "'65.3%"

$ws.Cells.Item(17, 18).Value = "'55.6%"
$ws.Cells.Item(17, 19).Value = "'61.5%"

$ws.Cells.Item(18, 19).Value = "'74.6%"
$ws.Cells.Item(19, 19).Value = "'75.9%"

# Restore clean (non quote-prefixed) formatting on every cell we just typed
# into with a leading apostrophe, copying it from an untouched style-4 cell.
$ws.Range("K9").Copy()
$ws.Range("L9").PasteSpecial(-4122)
$ws.Range("L10").PasteSpecial(-4122)
$ws.Range("R15").PasteSpecial(-4122)
$ws.Range("R16").PasteSpecial(-4122)
$ws.Range("R17").PasteSpecial(-4122)
$ws.Range("S15").PasteSpecial(-4122)
$ws.Range("S16").PasteSpecial(-4122)
$ws.Range("S17").PasteSpecial(-4122)
$ws.Range("S18").PasteSpecial(-4122)
$ws.Range("S19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Recorded / Pending counts (O/Q columns) and numeric student counts
$ws.Cells.Item(15, 15).Value = 15
$ws.Cells.Item(15, 17).Value = 12

$ws.Cells.Item(16, 15).Value = 14
$ws.Cells.Item(16, 17).Value = 12

$ws.Cells.Item(17, 15).Value = 15
$ws.Cells.Item(17, 17).Value = 12
